# Updated cryptos list with latest price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.135.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.302.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.53%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.38%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.663"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.67%  "
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.649.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.899"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.329.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.155.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.90%  "
$ws.Range("E20").Value = "  +7.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "238.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.01%  "
$ws.Range("E33").Value = "  +11.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0811"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +25.62%  "
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("E38").Value = "  +15.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.26%  "
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +21.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.212"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.36%  "
$ws.Range("E48").Value = "  +5.85%  "
$ws.Range("E49").Value = "  +5.53%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  +5.24%  "
